$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("B1").Value = "Di (mm)"
$ws.Range("C1").Value = "Xi exp"

# Extend column A style (bold, border, centered) to new rows 22:30, matching existing rows
$ws.Cells.Item(21, 1).Copy($ws.Range("A22:A30"))

# Update data for rows 2-30
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0.01309090777198838
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 5.201071428571429
$ws.Cells.Item(3, 3).Value = 0.02268346734018431
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 10.40214285714286
$ws.Cells.Item(4, 3).Value = 0.03762248338819871
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 15.60321428571429
$ws.Cells.Item(5, 3).Value = 0.05977083533143
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 20.80428571428572
$ws.Cells.Item(6, 3).Value = 0.0910311208835999
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 26.00535714285714
$ws.Cells.Item(7, 3).Value = 0.1330337733836711
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 31.20642857142857
$ws.Cells.Item(8, 3).Value = 0.1867608082790509
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 36.4075
$ws.Cells.Item(9, 3).Value = 0.2521854744786399
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 41.60857142857143
$ws.Cells.Item(10, 3).Value = 0.3280297364691219
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 46.80964285714286
$ws.Cells.Item(11, 3).Value = 0.4117318800502382
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 52.01071428571429
$ws.Cells.Item(12, 3).Value = 0.4996712014544753
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 57.21178571428572
$ws.Cells.Item(13, 3).Value = 0.5876265675321952
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 62.41285714285716
$ws.Cells.Item(14, 3).Value = 0.6713745342684518
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 67.61392857142859
$ws.Cells.Item(15, 3).Value = 0.7472880111307811
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 72.81500000000001
$ws.Cells.Item(16, 3).Value = 0.8127962811118726
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 78.01607142857144
$ws.Cells.Item(17, 3).Value = 0.8666116039870732
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 83.21714285714286
$ws.Cells.Item(18, 3).Value = 0.9086986318356469
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 88.41821428571428
$ws.Cells.Item(19, 3).Value = 0.9400331445060129
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 93.61928571428571
$ws.Cells.Item(20, 3).Value = 0.962242189540759
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 98.82035714285713
$ws.Cells.Item(21, 3).Value = 0.9772276098567154
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 104.0214285714286
$ws.Cells.Item(22, 3).Value = 0.9868534779720193
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 109.2225
$ws.Cells.Item(23, 3).Value = 0.9927397811900183
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 114.4235714285714
$ws.Cells.Item(24, 3).Value = 0.9961664875199779
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 119.6246428571428
$ws.Cells.Item(25, 3).Value = 0.9980655639659218
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 124.8257142857143
$ws.Cells.Item(26, 3).Value = 0.9990674978428111
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 130.0267857142857
$ws.Cells.Item(27, 3).Value = 0.9995707277667851
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 135.2278571428571
$ws.Cells.Item(28, 3).Value = 0.9998113436816461
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 140.4289285714285
$ws.Cells.Item(29, 3).Value = 0.9999208684709087
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 145.63
$ws.Cells.Item(30, 3).Value = 0.9999683287581669
